# Fix a typo in the Team ID on the cover page:
#   LTVIP2026TMIDS64787  ->  LTVIP2026TMIDS64789
# (last digit 7 -> 9; everything else - including the bold/28pt run
# formatting - stays the same).

$d = $word.ActiveDocument

$found = $d.Content.Find.Execute(
    "LTVIP2026TMIDS64787",   # FindText
    $true,                   # MatchCase
    $false,                  # MatchWholeWord
    $false,                  # MatchWildcards
    $false,                  # MatchSoundsLike
    $false,                  # MatchAllWordForms
    $true,                   # Forward
    1,                       # Wrap (wdFindContinue)
    $false,                  # Format
    "LTVIP2026TMIDS64789",   # ReplaceWith
    2                        # Replace (wdReplaceAll)
)

if (-not $found) {
    throw "Team ID text 'LTVIP2026TMIDS64787' was not found in the document."
}
